$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testSheet")

# Client name on the login/client row (I2, under "Client_Name" header) was
# changed from "Aron" to "Ajim".
$ws.Range("I2").Value = "Ajim"

# New "Add Client" verification message columns added next to the existing
# "Add_Reparation expected message" row.
$ws.Range("B7").Value = "Add_Client Message"
$ws.Range("B8").Value = "Client: Ajim"

# Selection moved to C7.
$ws.Range("C7").Select()
